# Update "想去人数" (interest count) figures in the "展览" and "全部类型"
# sheets to match the newly generated GitHub Pages data (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of cell -> new value, applied identically to both sheets.
$updates = @{
    "F5"  = 13049
    "F12" = 13752
    "F13" = 14310
    "F21" = 33
    "F23" = 111
    "F25" = 5375
    "F26" = 934
    "F28" = 305
    "F29" = 13
    "F30" = 22
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
